$d = $word.ActiveDocument

# Insert a brand-new, empty paragraph immediately before the document's
# current first paragraph, so it becomes the new first paragraph.
$d.Paragraphs(1).Range.InsertParagraphBefore()

# Give the new (now) first paragraph its text.
$newPara = $d.Paragraphs(1).Range
$newPara.Text = "IS-91 Khmelinin Andrey"

# Re-acquire the range (its extent grew when we set .Text) and format it:
# Times New Roman, bold, 14pt, with the same hanging/negative indent used
# throughout the rest of the document.
$newPara = $d.Paragraphs(1).Range
$newPara.Font.Name = "Times New Roman"
$newPara.Font.Bold = $true
$newPara.Font.BoldBi = $true
$newPara.Font.Size = 14
$newPara.ParagraphFormat.LeftIndent = -58.5
$newPara.ParagraphFormat.FirstLineIndent = 36
